# Update "想去人数" (interest count) figures in the 展览 and 全部类型 sheets
# to reflect newly generated output, per commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map: worksheet name -> hashtable of row -> new value for column F
$updates = @{
    "展览" = @{
        3  = 194
        4  = 580
        6  = 469
        9  = 2416
        10 = 140
        12 = 160
        13 = 1463
        19 = 179
        23 = 4
        24 = 115
        26 = 1512
        28 = 376
        29 = 317
    }
    "全部类型" = @{
        3  = 194
        4  = 580
        7  = 469
        10 = 2416
        11 = 140
        13 = 160
        14 = 1463
        20 = 179
        24 = 4
        25 = 115
        27 = 1512
        29 = 376
        30 = 317
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
